# Insert a new data row before row 302. This shifts the existing rows
# 302..377 down to 303..378 (carrying all of their original values with
# them), and leaves a blank row 302 ready to be populated with the new
# observation.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(302).Insert()

# The row that used to be 302 is now row 303, holding every original
# value untouched. Duplicate that row's contents into the freshly
# inserted row 302 so every column (Mercado ID, Mercado, Region, Codreg,
# Categoria ID, Categoria, Variedad, Calidad, prices, Unidad, Origen,
# Precio $/Kg, Kg o Unidades, Clasificacion, ...) starts out identical.
$ws.Range("A303:R303").Copy($ws.Range("A302:R302"))

# Now overwrite just the two cells that hold genuinely new data for this
# new record: the date (Fecha) and the Volumen.
$ws.Range("D302").Value = 44722
$ws.Range("J302").Value = 750
